# Lab Exam 03 grading — fill in "Points for grading" (column E) for the
# "Generic" and "Customer Class" rubric sections (rows 3-6 and 10-14),
# matching the full marks already entered in column D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 1
$ws.Range("E4").Value = 2
$ws.Range("E5").Value = 2
$ws.Range("E6").Value = 2

$ws.Range("E10").Value = 2
$ws.Range("E11").Value = 2
$ws.Range("E12").Value = 2
$ws.Range("E13").Value = 2
$ws.Range("E14").Value = 2

$ws.Range("E15").Select() | Out-Null
